$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.362.01'
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('D3').Value = '1.795.48'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '1.004'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '307.01'
$ws.Range('E6').Value = '  -1.20%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4519'
$ws.Range('E7').Value = '  -1.21%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3595'
$ws.Range('E8').Value = '  -2.26%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '46.29'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07082'
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8852'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '19.48'
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('D14').Value = '1.801.41'
$ws.Range('E14').Value = '  -1.06%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.283'
$ws.Range('E15').Value = '  -0.61%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.324'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '84.93'
$ws.Range('E17').Value = '  -2.16%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.006'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000008508'
$ws.Range('E19').Value = '  -2.16%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.004'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.24'
$ws.Range('E21').Value = '  -1.39%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').Value = '26.393.28'
$ws.Range('E22').Value = '  -1.86%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.971'
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '2.031.16'
$ws.Range('E24').Value = '  -0.67%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '10.54'
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.969'
$ws.Range('E26').Value = '  -1.79%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '151.33'
$ws.Range('E27').Value = '  +0.18%  '
$ws.Range('E28').Value = '  -1.84%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.026'
$ws.Range('E29').Value = '  +3.86%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '111.91'
$ws.Range('E30').Value = '  -1.38%  '
$ws.Range('E31').Value = '  -1.26%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.08685'
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.089'
$ws.Range('E33').Value = '  +2.39%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.769'
$ws.Range('E34').Value = '  +9.12%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.443'
$ws.Range('E35').Value = '  -0.54%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7225'
$ws.Range('E36').Value = '  -3.53%  '
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.003'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.066'
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01934'
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.852'
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.5067'
$ws.Range('E43').Value = '  +2.01%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '6.833'
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1517'
$ws.Range('E45').Value = '  -4.70%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '8.022'
$ws.Range('E46').Value = '  -3.15%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.004'
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.876'
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '100.58'
$ws.Range('E50').Value = '  -0.67%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.564'
$ws.Range('E51').Value = '  -2.79%  '
